$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 (SIQ_10): question reworded, asked-by and date updated
$ws.Range("B13").Value = "Amr"
$ws.Range("D13").Value = "What to do when alarm is reached if mode is in another mode?"
$ws.Range("G13").Value = "22/2/2020"

# Row 17 (SIQ_14): new question added, requirement id set, asked-by swapped, proposed-by cleared
$ws.Range("A17").Value = "PO3_DGW_CRS_F_02"
$ws.Range("B17").Value = "Mariam"
$ws.Range("D17").Value = "If user doesn't stop alarm manually, for how long does it keep ringing?"
$ws.Range("F17").ClearContents()

# Remove the last (blank) row of the table so it matches the new extent
$ws.Rows("38:38").Delete()
